# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Mon Jul 15 16:21:02 UTC 2024 with GitHub Actions"
#
# Column D ("Price") values are written through a temporary formula and then
# converted in-place to a literal value via Copy + PasteSpecial(values only).
# This keeps values such as "574.10", "0.999" or "1.00" stored as plain text
# (matching the original inlineStr cells) instead of being auto-coerced into
# numbers by a direct .Value assignment, and it avoids creating any new/unused
# cell styles (NumberFormat tricks leave orphaned style records behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link text updates (two rows swapped rank position) -----------------
$textUpdates = @(
    @{Cell="B32"; Value='NEARProtocol'},
    @{Cell="C32"; Value='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'},
    @{Cell="B33"; Value='RenderToken'},
    @{Cell="C33"; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'},
    @{Cell="B48"; Value='Bittensor'},
    @{Cell="C48"; Value='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'},
    @{Cell="B49"; Value='InjectiveProtocol'},
    @{Cell="C49"; Value='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'}
)
foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Price (column D) updates ----------------------------------------------------
$priceUpdates = @(
    @{Cell="D2"; Value='63.083.04'},
    @{Cell="D3"; Value='3.372.66'},
    @{Cell="D5"; Value='574.10'},
    @{Cell="D6"; Value='153.51'},
    @{Cell="D7"; Value='0.999'},
    @{Cell="D8"; Value='3.378.00'},
    @{Cell="D13"; Value='3.950.77'},
    @{Cell="D15"; Value='0.0000181'},
    @{Cell="D16"; Value='26.96'},
    @{Cell="D17"; Value='63.129.09'},
    @{Cell="D18"; Value='3.381.27'},
    @{Cell="D19"; Value='6.34'},
    @{Cell="D20"; Value='13.89'},
    @{Cell="D21"; Value='8.39'},
    @{Cell="D22"; Value='385.81'},
    @{Cell="D23"; Value='1.00'},
    @{Cell="D24"; Value='0.535'},
    @{Cell="D25"; Value='70.37'},
    @{Cell="D26"; Value='9.38'},
    @{Cell="D28"; Value='0.0₃0986'},
    @{Cell="D31"; Value='23.12'},
    @{Cell="D32"; Value='5.58'},
    @{Cell="D33"; Value='6.34'},
    @{Cell="D34"; Value='1.30'},
    @{Cell="D35"; Value='6.74'},
    @{Cell="D37"; Value='157.94'},
    @{Cell="D38"; Value='1.88'},
    @{Cell="D39"; Value='27.48'},
    @{Cell="D40"; Value='2.884.21'},
    @{Cell="D41"; Value='0.0744'},
    @{Cell="D42"; Value='0.0326'},
    @{Cell="D43"; Value='40.87'},
    @{Cell="D44"; Value='0.749'},
    @{Cell="D45"; Value='4.26'},
    @{Cell="D47"; Value='3.416.86'},
    @{Cell="D48"; Value='301.90'},
    @{Cell="D49"; Value='21.94'}
)
foreach ($u in $priceUpdates) {
    $r = $ws.Range($u.Cell)
    $r.Formula = '="' + $u.Value + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

# --- Volume(1h) percentage (column E) updates --------------------------------------
$volumeUpdates = @(
    @{Cell="E2"; Value='  +5.05%  '},
    @{Cell="E3"; Value='  +5.64%  '},
    @{Cell="E4"; Value='  -0.01%  '},
    @{Cell="E5"; Value='  +7.01%  '},
    @{Cell="E6"; Value='  +5.41%  '},
    @{Cell="E7"; Value='  -0.17%  '},
    @{Cell="E8"; Value='  +5.57%  '},
    @{Cell="E9"; Value='  -0.61%  '},
    @{Cell="E10"; Value='  +1.59%  '},
    @{Cell="E11"; Value='  +5.50%  '},
    @{Cell="E12"; Value='  +1.02%  '},
    @{Cell="E13"; Value='  +5.48%  '},
    @{Cell="E14"; Value='  +0.16%  '},
    @{Cell="E15"; Value='  +4.51%  '},
    @{Cell="E16"; Value='  +4.01%  '},
    @{Cell="E17"; Value='  +5.06%  '},
    @{Cell="E18"; Value='  +5.78%  '},
    @{Cell="E19"; Value='  +1.01%  '},
    @{Cell="E20"; Value='  +4.60%  '},
    @{Cell="E21"; Value='  +2.03%  '},
    @{Cell="E22"; Value='  +4.49%  '},
    @{Cell="E23"; Value='  -0.06%  '},
    @{Cell="E24"; Value='  +2.28%  '},
    @{Cell="E25"; Value='  +1.36%  '},
    @{Cell="E26"; Value='  +8.43%  '},
    @{Cell="E27"; Value='  +5.91%  '},
    @{Cell="E28"; Value='  +12.74%  '},
    @{Cell="E29"; Value='  +1.21%  '},
    @{Cell="E30"; Value='  +6.90%  '},
    @{Cell="E31"; Value='  +3.05%  '},
    @{Cell="E32"; Value='  +5.68%  '},
    @{Cell="E33"; Value='  +3.88%  '},
    @{Cell="E34"; Value='  +9.07%  '},
    @{Cell="E35"; Value='  +2.65%  '},
    @{Cell="E36"; Value='  +9.29%  '},
    @{Cell="E37"; Value='  +0.79%  '},
    @{Cell="E38"; Value='  +12.27%  '},
    @{Cell="E39"; Value='  +5.54%  '},
    @{Cell="E40"; Value='  +2.14%  '},
    @{Cell="E41"; Value='  +5.49%  '},
    @{Cell="E42"; Value='  +8.12%  '},
    @{Cell="E43"; Value='  +2.60%  '},
    @{Cell="E44"; Value='  +4.33%  '},
    @{Cell="E45"; Value='  +0.43%  '},
    @{Cell="E46"; Value='  +5.91%  '},
    @{Cell="E47"; Value='  +5.66%  '},
    @{Cell="E48"; Value='  +14.11%  '},
    @{Cell="E49"; Value='  +5.97%  '},
    @{Cell="E50"; Value='  -1.51%  '},
    @{Cell="E51"; Value='  +1.89%  '}
)
foreach ($u in $volumeUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

